$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'35.381.81"
$ws.Cells.Item(2, 5).Value = "  +0.49%  "
$ws.Cells.Item(3, 4).Value = "'1.910.58"
$ws.Cells.Item(3, 5).Value = "  +2.75%  "
$ws.Cells.Item(4, 5).Value = "  -0.50%  "
$ws.Cells.Item(5, 4).Value = "'246.34"
$ws.Cells.Item(5, 5).Value = "  +2.87%  "
$ws.Cells.Item(6, 5).Value = "  +5.98%  "
$ws.Cells.Item(7, 5).Value = "  -0.51%  "
$ws.Cells.Item(8, 4).Value = "'41.37"
$ws.Cells.Item(8, 5).Value = "  -1.80%  "
$ws.Cells.Item(9, 5).Value = "  +6.19%  "
$ws.Cells.Item(10, 4).Value = "'52.78"
$ws.Cells.Item(10, 5).Value = "  +12.53%  "
$ws.Cells.Item(11, 4).Value = "'0.0720"
$ws.Cells.Item(11, 5).Value = "  +3.88%  "
$ws.Cells.Item(12, 5).Value = "  +0.29%  "
$ws.Cells.Item(13, 4).Value = "'2.187.13"
$ws.Cells.Item(13, 5).Value = "  +2.73%  "
$ws.Cells.Item(14, 4).Value = "'12.12"
$ws.Cells.Item(14, 5).Value = "  +5.48%  "
$ws.Cells.Item(15, 4).Value = "'0.702"
$ws.Cells.Item(15, 5).Value = "  +3.84%  "
$ws.Cells.Item(16, 4).Value = "'1.903.22"
$ws.Cells.Item(16, 5).Value = "  +2.22%  "
$ws.Cells.Item(17, 5).Value = "  +3.16%  "
$ws.Cells.Item(18, 4).Value = "'35.356.22"
$ws.Cells.Item(18, 5).Value = "  +0.49%  "
$ws.Cells.Item(19, 4).Value = "'72.58"
$ws.Cells.Item(19, 5).Value = "  +3.83%  "
$ws.Cells.Item(20, 4).Value = "'0.0₃0826"
$ws.Cells.Item(20, 5).Value = "  +3.56%  "
$ws.Cells.Item(21, 4).Value = "'239.94"
$ws.Cells.Item(21, 5).Value = "  -0.27%  "
$ws.Cells.Item(22, 5).Value = "  +2.62%  "
$ws.Cells.Item(23, 5).Value = "  +2.08%  "
$ws.Cells.Item(24, 5).Value = "  -0.42%  "
$ws.Cells.Item(25, 4).Value = "'2.30"
$ws.Cells.Item(25, 5).Value = "  +1.19%  "
$ws.Cells.Item(26, 4).Value = "'2.36"
$ws.Cells.Item(26, 5).Value = "  +23.27%  "
$ws.Cells.Item(27, 4).Value = "'169.64"
$ws.Cells.Item(27, 5).Value = "  +0.56%  "
$ws.Cells.Item(28, 4).Value = "'8.48"
$ws.Cells.Item(28, 5).Value = "  +5.92%  "
$ws.Cells.Item(29, 5).Value = "  +4.82%  "
$ws.Cells.Item(31, 4).Value = "'4.17"
$ws.Cells.Item(31, 5).Value = "  +4.04%  "
$ws.Cells.Item(32, 4).Value = "'0.0567"
$ws.Cells.Item(32, 5).Value = "  +1.42%  "
$ws.Cells.Item(33, 4).Value = "'0.943"
$ws.Cells.Item(33, 5).Value = "  +15.74%  "
$ws.Cells.Item(34, 5).Value = "  -0.45%  "
$ws.Cells.Item(35, 5).Value = "  +2.56%  "
$ws.Cells.Item(36, 5).Value = "  -4.00%  "
$ws.Cells.Item(37, 4).Value = "'2.05"
$ws.Cells.Item(37, 5).Value = "  +0.54%  "
$ws.Cells.Item(38, 5).Value = "  +1.42%  "
$ws.Cells.Item(39, 5).Value = "  +1.97%  "
$ws.Cells.Item(40, 5).Value = "  +3.87%  "
$ws.Cells.Item(41, 4).Value = "'16.26"
$ws.Cells.Item(41, 5).Value = "  +8.70%  "
$ws.Cells.Item(42, 4).Value = "'0.0639"
$ws.Cells.Item(42, 5).Value = "  +8.53%  "
$ws.Cells.Item(43, 4).Value = "'90.22"
$ws.Cells.Item(43, 5).Value = "  +0.33%  "
$ws.Cells.Item(44, 4).Value = "'1.340.38"
$ws.Cells.Item(45, 5).Value = "  +3.20%  "
$ws.Cells.Item(46, 4).Value = "'48.00"
$ws.Cells.Item(46, 5).Value = "  +37.99%  "
$ws.Cells.Item(47, 5).Value = "  +1.86%  "
$ws.Cells.Item(48, 4).Value = "'2.40"
$ws.Cells.Item(48, 5).Value = "  -0.78%  "
$ws.Cells.Item(49, 5).Value = "  -0.26%  "
$ws.Cells.Item(50, 4).Value = "'2.094.76"
$ws.Cells.Item(50, 5).Value = "  +2.35%  "
$ws.Cells.Item(51, 4).Value = "'0.0704"
$ws.Cells.Item(51, 5).Value = "  +3.59%  "
